$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.52879766666667
$ws.Range("H2").Value = 40.586393
$ws.Range("I2").Value = 0.8596824606989164
$ws.Range("J2").Value = 0.8596824606989165
$ws.Range("M2").Value = 60.538204
$ws.Range("N2").Value = 181.614612
$ws.Range("O2").Value = 0.6123615450168176
$ws.Range("P2").Value = 0.6123615450168176
$ws.Range("Q2").Value = 819.0091130193907
$ws.Range("R2").Value = 7371.082017174516
$ws.Range("S2").Value = 0.526436479857448
$ws.Range("T2").Value = 0.5264364798574481
$ws.Range("G3").Value = 13.52879766666667
$ws.Range("H3").Value = 40.586393
$ws.Range("I3").Value = 0.8596824606989164
$ws.Range("J3").Value = 0.8596824606989165
$ws.Range("O3").Value = 0.1096681363892149
$ws.Range("P3").Value = 0.1096681363892149
$ws.Range("Q3").Value = 146.6767530416262
$ws.Range("R3").Value = 1320.090777374636
$ws.Range("S3").Value = 0.09427977335134462
$ws.Range("T3").Value = 0.09427977335134465
$ws.Range("G4").Value = 13.52879766666667
$ws.Range("H4").Value = 40.586393
$ws.Range("I4").Value = 0.8596824606989164
$ws.Range("J4").Value = 0.8596824606989165
$ws.Range("M4").Value = 8.850437666666666
$ws.Range("N4").Value = 26.551313
$ws.Range("O4").Value = 0.08952475173586316
$ws.Range("P4").Value = 0.08952475173586316
$ws.Range("Q4").Value = 119.7357804537788
$ws.Range("R4").Value = 1077.622024084009
$ws.Range("S4").Value = 0.07696285886574643
$ws.Range("T4").Value = 0.07696285886574644
$ws.Range("G5").Value = 13.52879766666667
$ws.Range("H5").Value = 40.586393
$ws.Range("I5").Value = 0.8596824606989164
$ws.Range("J5").Value = 0.8596824606989165
$ws.Range("M5").Value = 1.757142
$ws.Range("N5").Value = 5.271426
$ws.Range("O5").Value = 0.01777400251143792
$ws.Range("P5").Value = 0.01777400251143792
$ws.Range("Q5").Value = 23.772018589602
$ws.Range("R5").Value = 213.948167306418
$ws.Range("S5").Value = 0.01527999821550167
$ws.Range("T5").Value = 0.01527999821550167
$ws.Range("G6").Value = 13.52879766666667
$ws.Range("H6").Value = 40.586393
$ws.Range("I6").Value = 0.8596824606989164
$ws.Range("J6").Value = 0.8596824606989165
$ws.Range("M6").Value = 16.87263033333333
$ws.Range("N6").Value = 50.617891
$ws.Range("O6").Value = 0.1706715643466665
$ws.Range("P6").Value = 0.1706715643466665
$ws.Range("Q6").Value = 228.2664018841292
$ws.Range("R6").Value = 2054.397616957163
$ws.Range("S6").Value = 0.1467233504088757
$ws.Range("T6").Value = 0.1467233504088757
$ws.Range("I7").Value = 0.008772717629348061
$ws.Range("J7").Value = 0.008772717629348063
$ws.Range("M7").Value = 60.538204
$ws.Range("N7").Value = 181.614612
$ws.Range("O7").Value = 0.6123615450168176
$ws.Range("P7").Value = 0.6123615450168176
$ws.Range("Q7").Value = 8.357662291423999
$ws.Range("R7").Value = 75.218960622816
$ws.Range("S7").Value = 0.005372074921503852
$ws.Range("T7").Value = 0.005372074921503853
$ws.Range("I8").Value = 0.008772717629348061
$ws.Range("J8").Value = 0.008772717629348063
$ws.Range("O8").Value = 0.1096681363892149
$ws.Range("P8").Value = 0.1096681363892149
$ws.Range("S8").Value = 0.0009620875934794129
$ws.Range("T8").Value = 0.0009620875934794132
$ws.Range("I9").Value = 0.008772717629348061
$ws.Range("J9").Value = 0.008772717629348063
$ws.Range("M9").Value = 8.850437666666666
$ws.Range("N9").Value = 26.551313
$ws.Range("O9").Value = 0.08952475173586316
$ws.Range("P9").Value = 0.08952475173586316
$ws.Range("Q9").Value = 1.221856022509333
$ws.Range("R9").Value = 10.996704202584
$ws.Range("S9").Value = 0.0007853753678162151
$ws.Range("T9").Value = 0.0007853753678162154
$ws.Range("I10").Value = 0.008772717629348061
$ws.Range("J10").Value = 0.008772717629348063
$ws.Range("M10").Value = 1.757142
$ws.Range("N10").Value = 5.271426
$ws.Range("O10").Value = 0.01777400251143792
$ws.Range("P10").Value = 0.01777400251143792
$ws.Range("Q10").Value = 0.242583995952
$ws.Range("R10").Value = 2.183255963568
$ws.Range("S10").Value = 0.0001559263051761681
$ws.Range("T10").Value = 0.0001559263051761682
$ws.Range("I11").Value = 0.008772717629348061
$ws.Range("J11").Value = 0.008772717629348063
$ws.Range("M11").Value = 16.87263033333333
$ws.Range("N11").Value = 50.617891
$ws.Range("O11").Value = 0.1706715643466665
$ws.Range("P11").Value = 0.1706715643466665
$ws.Range("Q11").Value = 2.329367853298666
$ws.Range("R11").Value = 20.964310679688
$ws.Range("S11").Value = 0.001497253441372413
$ws.Range("T11").Value = 0.001497253441372413
$ws.Range("G12").Value = 0.05669233333333334
$ws.Range("H12").Value = 0.170077
$ws.Range("I12").Value = 0.003602493423554283
$ws.Range("J12").Value = 0.003602493423554284
$ws.Range("M12").Value = 60.538204
$ws.Range("N12").Value = 181.614612
$ws.Range("O12").Value = 0.6123615450168176
$ws.Range("P12").Value = 0.6123615450168176
$ws.Range("Q12").Value = 3.432052040569334
$ws.Range("R12").Value = 30.888468365124
$ws.Range("S12").Value = 0.002206028438760626
$ws.Range("T12").Value = 0.002206028438760626
$ws.Range("G13").Value = 0.05669233333333334
$ws.Range("H13").Value = 0.170077
$ws.Range("I13").Value = 0.003602493423554283
$ws.Range("J13").Value = 0.003602493423554284
$ws.Range("O13").Value = 0.1096681363892149
$ws.Range("P13").Value = 0.1096681363892149
$ws.Range("Q13").Value = 0.6146479222004445
$ws.Range("R13").Value = 5.531831299804001
$ws.Range("S13").Value = 0.0003950787401156007
$ws.Range("T13").Value = 0.0003950787401156009
$ws.Range("G14").Value = 0.05669233333333334
$ws.Range("H14").Value = 0.170077
$ws.Range("I14").Value = 0.003602493423554283
$ws.Range("J14").Value = 0.003602493423554284
$ws.Range("M14").Value = 8.850437666666666
$ws.Range("N14").Value = 26.551313
$ws.Range("O14").Value = 0.08952475173586316
$ws.Range("P14").Value = 0.08952475173586316
$ws.Range("Q14").Value = 0.5017519623445555
$ws.Range("R14").Value = 4.515767661101
$ws.Range("S14").Value = 0.000322512329373777
$ws.Range("T14").Value = 0.000322512329373777
$ws.Range("G15").Value = 0.05669233333333334
$ws.Range("H15").Value = 0.170077
$ws.Range("I15").Value = 0.003602493423554283
$ws.Range("J15").Value = 0.003602493423554284
$ws.Range("M15").Value = 1.757142
$ws.Range("N15").Value = 5.271426
$ws.Range("O15").Value = 0.01777400251143792
$ws.Range("P15").Value = 0.01777400251143792
$ws.Range("Q15").Value = 0.09961647997800001
$ws.Range("R15").Value = 0.896548319802
$ws.Range("S15").Value = 0.00006403072715769241
$ws.Range("T15").Value = 0.00006403072715769243
$ws.Range("G16").Value = 0.05669233333333334
$ws.Range("H16").Value = 0.170077
$ws.Range("I16").Value = 0.003602493423554283
$ws.Range("J16").Value = 0.003602493423554284
$ws.Range("M16").Value = 16.87263033333333
$ws.Range("N16").Value = 50.617891
$ws.Range("O16").Value = 0.1706715643466665
$ws.Range("P16").Value = 0.1706715643466665
$ws.Range("Q16").Value = 0.9565487830674445
$ws.Range("R16").Value = 8.608939047607
$ws.Range("S16").Value = 0.0006148431881465876
$ws.Range("T16").Value = 0.0006148431881465878
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.013424666666667
$ws.Range("H17").Value = 6.040274
$ws.Range("I17").Value = 0.1279423282481813
$ws.Range("J17").Value = 0.1279423282481813
$ws.Range("M17").Value = 60.538204
$ws.Range("N17").Value = 181.614612
$ws.Range("O17").Value = 0.6123615450168176
$ws.Range("P17").Value = 0.6123615450168176
$ws.Range("Q17").Value = 121.8891132092987
$ws.Range("R17").Value = 1097.002018883688
$ws.Range("S17").Value = 0.07834696179910511
$ws.Range("T17").Value = 0.07834696179910514
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 2.013424666666667
$ws.Range("H18").Value = 6.040274
$ws.Range("I18").Value = 0.1279423282481813
$ws.Range("J18").Value = 0.1279423282481813
$ws.Range("O18").Value = 0.1096681363892149
$ws.Range("P18").Value = 0.1096681363892149
$ws.Range("Q18").Value = 21.82918245042756
$ws.Range("R18").Value = 196.462642053848
$ws.Range("S18").Value = 0.01403119670427524
$ws.Range("T18").Value = 0.01403119670427525
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 2.013424666666667
$ws.Range("H19").Value = 6.040274
$ws.Range("I19").Value = 0.1279423282481813
$ws.Range("J19").Value = 0.1279423282481813
$ws.Range("M19").Value = 8.850437666666666
$ws.Range("N19").Value = 26.551313
$ws.Range("O19").Value = 0.08952475173586316
$ws.Range("P19").Value = 0.08952475173586316
$ws.Range("Q19").Value = 17.81968950886245
$ws.Range("R19").Value = 160.377205579762
$ws.Range("S19").Value = 0.01145400517292674
$ws.Range("T19").Value = 0.01145400517292674
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 2.013424666666667
$ws.Range("H20").Value = 6.040274
$ws.Range("I20").Value = 0.1279423282481813
$ws.Range("J20").Value = 0.1279423282481813
$ws.Range("M20").Value = 1.757142
$ws.Range("N20").Value = 5.271426
$ws.Range("O20").Value = 0.01777400251143792
$ws.Range("P20").Value = 0.01777400251143792
$ws.Range("Q20").Value = 3.537873045636
$ws.Range("R20").Value = 31.840857410724
$ws.Range("S20").Value = 0.002274047263602388
$ws.Range("T20").Value = 0.002274047263602389
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 2.013424666666667
$ws.Range("H21").Value = 6.040274
$ws.Range("I21").Value = 0.1279423282481813
$ws.Range("J21").Value = 0.1279423282481813
$ws.Range("M21").Value = 16.87263033333333
$ws.Range("N21").Value = 50.617891
$ws.Range("O21").Value = 0.1706715643466665
$ws.Range("P21").Value = 0.1706715643466665
$ws.Range("Q21").Value = 33.97177010468156
$ws.Range("R21").Value = 305.745930942134
$ws.Range("S21").Value = 0.02183611730827179
$ws.Range("T21").Value = 0.0218361173082718
